$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A1")
$rng.Font.Name = "Times New Roman"
$rng.Font.Size = 12
$rng.Borders.LineStyle = "Continuous"
$rng.HorizontalAlignment = "General"
$rng.VerticalAlignment = "Bottom"
$rng.WrapText = $false
$rng.Value = "Test"
